$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is reshaped from two columns (searchText/result pairs with a
# jenkins/pine lookup and numeric results) into a single vertical column of
# status values. Drop column B's old "result" header and numeric values...
$ws.Range("B1:B3").ClearContents()

# ...and restack column A as status / available / pending / sold.
$ws.Range("A1").Value = "status"
$ws.Range("A2").Value = "available"
$ws.Range("A3").Value = "pending"
$ws.Range("A4").Value = "sold"

# Selection moves on to the next empty row, matching the saved file's cursor.
$ws.Range("A5").Select()
